$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set in_service (column E) to TRUE for rows 10 through 14
$ws.Range("E10:E14").Value = $true
